$d = $word.ActiveDocument

# 1) Expand the DIDs clause to mention labeled / property graph statements
$d.Content.Find.Execute(
    "transforms, saga / zippers",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "transforms, labeled / property graph statements: saga / zippers",
    2
)

# 2) Rework the Zippers clause: colon -> comma, and add "/ labeled property graphs"
$d.Content.Find.Execute(
    "mutable chain branches: dimensional contexts).",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "mutable chain branches, dimensional contexts / labeled property graphs).",
    2
)
